$d = $word.ActiveDocument

# Curly quote characters used throughout this document.
$ldq = [char]0x201C   # “
$rdq = [char]0x201D   # ”

# Locate the paragraph ending in "...we exit from the iteration." (unique anchor).
$range = $d.Content
$range.Find.ClearFormatting()
$anchor = "we exit from the iteration."
$found = $range.Find.Execute($anchor, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    # Collapse to the end of the matched text (end of the sentence / paragraph).
    $range.Collapse(0)

    $newText = " The " + $ldq + "iostat_end" + $rdq + " constant is defined in the intrinsic module " + $ldq + "iso_fotran_env" + $rdq + "."

    $range.InsertAfter($newText)

    # Force the newly inserted text into its own run (distinct from the preceding
    # run) by toggling a character property on it and then clearing it back to
    # "not set" (wdUndefined), mirroring the paragraph's original (unformatted) runs.
    $range.Font.Bold = 1
    $range.Font.Bold = 9999999
}
